# Regenerate s_val data to filter save games.
# Updates columns B (TB), C (d2S), D (K), E (IP) and the derived sum
# column G for each data row (2-12). Column A (date) and F (Win) are
# unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ B = 0.001754667048134761;  C = 0.05231270169004087;  D = 0.7127328510149897;  E = 246.9852506941017 }
    3  = @{ B = 3.182878228561681;     C = 1.65323645889881;     D = 0.7127328510149897;  E = 0.4998867070740569 }
    4  = @{ B = 3.182878228561681;     C = 1.65323645889881;     D = 0.7127328510149897;  E = 0.4998867070740569 }
    5  = @{ B = 0.7287194209349384;    C = 0.3375848360084654;   D = 0.7127328510149897;  E = 0.4998867070740569 }
    6  = @{ B = 3.182878228561681;     C = 1.65323645889881;     D = 3.082599426703578;   E = 0.4998867070740569 }
    7  = @{ B = 3.182878228561681;     C = 1.65323645889881;     D = 0.1529057820181812;  E = 0.4998867070740569 }
    8  = @{ B = 3.182878228561681;     C = 9.226618575922256;    D = 16.98373111632243;   E = 6.48142807727062 }
    9  = @{ B = 3.182878228561681;     C = 1.65323645889881;     D = 0.1529057820181812;  E = 0.4998867070740569 }
    10 = @{ B = 1.505614041169197;     C = 1.65323645889881;     D = 0.7127328510149897;  E = 0.4998867070740569 }
    11 = @{ B = 0.3464964993005633;    C = 9.226618575922256;    D = 0.7127328510149897;  E = 6.48142807727062 }
    12 = @{ B = 0.00006486019690155054; C = 0.05231270169004087;  D = 0.1529057820181812;  E = 6.48142807727062 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals.B
    $ws.Cells.Item($row, 3).Value = $vals.C
    $ws.Cells.Item($row, 4).Value = $vals.D
    $ws.Cells.Item($row, 5).Value = $vals.E
    $ws.Cells.Item($row, 7).Value = $vals.B + $vals.C + $vals.D + $vals.E
}
